$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '43.933.64'
$cell.Style = "Normal"

$ws.Range("E2").Value = '  -0.16%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.370.59'
$cell.Style = "Normal"

$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.05%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '0.675'
$cell.Style = "Normal"

$ws.Range("E5").Value = '  -0.68%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '240.66'
$cell.Style = "Normal"

$ws.Range("E6").Value = '  +0.55%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '74.29'
$cell.Style = "Normal"

$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").Value = '  +0.06%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.618'
$cell.Style = "Normal"

$ws.Range("E9").Value = '  +3.32%  '
$ws.Range("E10").Value = '  +2.31%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '60.55'
$cell.Style = "Normal"

$ws.Range("E11").Value = '  +5.66%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '37.90'
$cell.Style = "Normal"

$ws.Range("E12").Value = '  +16.46%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '7.33'
$cell.Style = "Normal"

$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("E14").Value = '  +1.09%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '16.43'
$cell.Style = "Normal"

$ws.Range("E15").Value = '  -1.00%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.924'
$cell.Style = "Normal"

$ws.Range("E16").Value = '  +2.84%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.371.56'
$cell.Style = "Normal"

$ws.Range("E17").Value = '  +0.56%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '43.928.28'
$cell.Style = "Normal"

$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("E19").Value = '  +2.04%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '78.15'
$cell.Style = "Normal"

$ws.Range("E20").Value = '  +1.62%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.59'
$cell.Style = "Normal"

$ws.Range("E21").Value = '  -2.31%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '254.51'
$cell.Style = "Normal"

$ws.Range("E22").Value = '  -1.55%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +3.36%  '
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("E26").Value = '  +0.45%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '10.57'
$cell.Style = "Normal"

$ws.Range("E27").Value = '  -1.72%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '22.45'
$cell.Style = "Normal"

$ws.Range("E29").Value = '  -0.84%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '176.13'
$cell.Style = "Normal"

$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  -1.23%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0755'
$cell.Style = "Normal"

$ws.Range("E33").Value = '  -1.17%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '5.42'
$cell.Style = "Normal"

$ws.Range("E34").Value = '  -0.58%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '5.11'
$cell.Style = "Normal"

$ws.Range("E35").Value = '  -2.09%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '3.83'
$cell.Style = "Normal"

$ws.Range("E36").Value = '  +1.94%  '
$ws.Range("E37").Value = '  +5.28%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.43'
$cell.Style = "Normal"

$ws.Range("E38").Value = '  +3.14%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.0280'
$cell.Style = "Normal"

$ws.Range("E39").Value = '  +0.82%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '5.45'
$cell.Style = "Normal"

$ws.Range("E40").Value = '  +15.35%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '20.48'
$cell.Style = "Normal"

$ws.Range("E41").Value = '  +8.10%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '65.12'
$cell.Style = "Normal"

$ws.Range("E42").Value = '  +12.03%  '
$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.108'
$cell.Style = "Normal"

$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.204'
$cell.Style = "Normal"

$ws.Range("E44").Value = '  -1.23%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '9.10'
$cell.Style = "Normal"

$ws.Range("E45").Value = '  +1.05%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.57'
$cell.Style = "Normal"

$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.25'
$cell.Style = "Normal"

$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"

$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  -0.98%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '98.68'
$cell.Style = "Normal"

$ws.Range("E50").Value = '  -1.35%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '4.42'
$cell.Style = "Normal"

$ws.Range("E51").Value = '  +15.96%  '
